# "Generate Report for Handback"
#
# The handback-status report previously marked
# 76942e8c-7103-48cf-abd1-2a80b8dfde03.md as "Ready for handoff". A new
# report run discovered the handback transform failed for that file (the
# handback archive's file name did not match the expected handoff file
# name), so the status is updated everywhere it is shown, and the
# per-language "Error Detail" column is populated with the failure
# reason. The Error Detail column is also widened so the message is
# readable.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-language status for the failed file (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: Status + Error Detail for the failed file (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"
$wsZh.Range("P3").Value = "Handback file name: w2e5wzyq.m0v is different with handoff file name: 76942e8c-7103-48cf-abd1-2a80b8dfde03.9005920b1d370fc5fcd3736951a147e30654efcf.zh-cn."
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Status + Error Detail for the failed file (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"
$wsDe.Range("P3").Value = "Handback file name: w2e5wzyq.m0v is different with handoff file name: 76942e8c-7103-48cf-abd1-2a80b8dfde03.9005920b1d370fc5fcd3736951a147e30654efcf.de-de."
$wsDe.Columns.Item(16).ColumnWidth = 39.17
